# Add team record (Wins/Losses/Ties) columns to the data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (AC1) onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record values for every data row (2 through 40)
$lastRow = 40
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95
    $ws.Cells.Item($r, 31).Value = 67
    $ws.Cells.Item($r, 32).Value = 0
}
